# Data update using git
# Applies updated enrollment/payment counts to the "Resumo Inscricoes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("E5").Value = 165
$ws.Range("F5").Value = 113
$ws.Range("H5").Value = 124

# Row 10
$ws.Range("E10").Value = 704
$ws.Range("F10").Value = 394
$ws.Range("H10").Value = 489

# Row 11
$ws.Range("E11").Value = 466

# Row 12
$ws.Range("E12").Value = 706

# Row 13
$ws.Range("E13").Value = 168
$ws.Range("F13").Value = 94
$ws.Range("H13").Value = 128

# Row 15
$ws.Range("F15").Value = 96
$ws.Range("H15").Value = 147

# Row 25
$ws.Range("E25").Value = 334

# Row 27
$ws.Range("F27").Value = 212
$ws.Range("H27").Value = 294

# Row 30
$ws.Range("E30").Value = 255
$ws.Range("F30").Value = 159
$ws.Range("H30").Value = 211

# Row 33
$ws.Range("E33").Value = 332

# Row 34
$ws.Range("E34").Value = 252
$ws.Range("F34").Value = 178
$ws.Range("H34").Value = 216

# Row 42
$ws.Range("E42").Value = 467

# Row 43
$ws.Range("E43").Value = 143

# Row 44
$ws.Range("E44").Value = 374

# Row 45
$ws.Range("E45").Value = 182
$ws.Range("F45").Value = 101
$ws.Range("H45").Value = 140

# Row 46
$ws.Range("E46").Value = 392

# Row 47
$ws.Range("E47").Value = 544
$ws.Range("F47").Value = 307
$ws.Range("H47").Value = 399

# Row 50
$ws.Range("E50").Value = 286
$ws.Range("F50").Value = 158
$ws.Range("H50").Value = 231
